# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# This updates column G ("K") on Sheet1 with the newly computed values,
# leaving all other columns untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = [ordered]@{
    2 = 0
    3 = 1
    4 = 0
    5 = 1
    6 = 1
    7 = 3
    8 = 2
    9 = 1
    10 = 0
    11 = 2
    12 = 0
    13 = 2
    14 = 0
    15 = 0
    16 = 2
    17 = 2
    18 = 0
    19 = 1
    20 = 1
    21 = 1
    22 = 1
    23 = 1
    24 = 1
    25 = 0
    26 = 1
    27 = 1
    28 = 1
    29 = 1
    30 = 0
    31 = 1
    32 = 0
    33 = 1
    34 = 2
    35 = 1
    36 = 1
    37 = 2
    38 = 1
    39 = 0
    40 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
